# Fixed most imports to search by id instead of name.
#
# Insert a new "id" column at the front of the Passive Skills sheet, populate
# it with sequential numeric ids (1, 2, 3, ...), and repoint the
# parent_skill_id column from the parent's NAME (a string lookup) to the
# parent's new numeric id. This lets people rename skills in the sheet
# without breaking the parent/child relationship used by the importer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Passive Skills")

$lastRow = 10

# Snapshot the current (pre-shift) skill name -> new numeric id mapping,
# and the current (pre-shift) parent_skill_id (column G) name per row --
# both need to be read before the column insert shifts everything right.
$nameToId = @{}
$parentNameByRow = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $nameToId[$ws.Range("A" + $r).Value2] = ($r - 1)
    $parentNameByRow[$r] = $ws.Range("G" + $r).Value2
}

# Insert a new blank column before column A; this shifts the old A:J to B:K.
$ws.Columns("A:A").Insert()
$ws.Columns("A").ColumnWidth = 2.17

# Header for the new id column.
$ws.Range("A1").Value = "id"

# Fill in the sequential numeric id values, and fix up parent_skill_id
# (now column H, was column G) to hold the parent's numeric id instead of
# its name.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("A" + $r).Value = $r - 1

    $parentName = $parentNameByRow[$r]
    if ($parentName -ne $null -and $parentName -ne "") {
        $ws.Range("H" + $r).Value = $nameToId[$parentName]
    }
}

[void]$ws.Range("A1").Select()
